$wb = $excel.ActiveWorkbook

# --- Locate the currently-last week sheet ("Semana 10") ---
$semana10 = $wb.Worksheets.Item("Semana 10")

# --- 1) Duplicate it to create next week's sheet, placed right after it ---
#     (mirrors Excel's "Move or Copy... > Create a copy" placed After Semana 10)
$semana10.Copy($null, $semana10)
$semana11 = $wb.ActiveSheet
$semana11.Name = "Semana 11"

# The fresh copy still holds Semana 10's OLD entries for row 4/5 (Martes/Miercoles),
# so fill it in with the new week's real time entries (Lunes & Martes columns).
$semana11.Range("A4").Value = 0.33333333333333331
$semana11.Range("B4").Value = 0.51666666666666672
$semana11.Range("D4").Value = 0.39166666666666666
$semana11.Range("E4").Value = 0.44791666666666669

$semana11.Range("A5").Value = 0.53749999999999998
$semana11.Range("B5").Value = 0.74513888888888891
$semana11.Range("D5").Value = 0.52708333333333335
$semana11.Range("E5").Value = 0.6875

$semana11.Range("A6").Value = 0.75555555555555554
$semana11.Range("B6").Value = 0.84652777777777777

# Put the cursor where it was left after data entry on the new sheet
$semana11.Range("E14").Select()

# --- 2) Go back and finish filling in "Semana 10" (Miercoles/Jueves/Viernes cols) ---
$semana10.Range("G4").Value = 0.33333333333333331
$semana10.Range("H4").Value = 0.58333333333333337

$semana10.Range("J4").Value = 0.60416666666666663
$semana10.Range("K4").Value = 0.72916666666666663

$semana10.Range("M4").Value = 0.39583333333333331
$semana10.Range("N4").Value = 0.5

$semana10.Range("M5").Value = 0.54166666666666663
$semana10.Range("N5").Value = 0.70833333333333337

# Select-all on Semana 10 (as left after reviewing the completed week)
$semana10.Cells.Select()

# --- 3) Make the new week the active tab, matching the end-of-session state ---
$semana11.Activate()
